$d = $word.ActiveDocument

# The title paragraph (first paragraph) currently reads:
#   "Depth First Search Pseudo-Code"
# split across two underlined runs: "Depth First Search" and " Pseudo-Code".
# We need it to read "Depth First Search Pseudo Code" (hyphen -> space) with
# the trailing "Code" living in its own run, e.g.:
#   "Depth First Search" | " Pseudo " | "Code"
$titlePara = $d.Paragraphs.Item(1)
$paraStart = $titlePara.Range.Start
$paraEnd = $titlePara.Range.End

# Step 1: Replace the hyphen in "Pseudo-Code" with a space so the text reads
# "... Pseudo Code". A plain character replacement like this can cause the
# engine to coalesce adjacent same-formatted runs into one run; that's fine,
# we fix the run boundaries in step 2.
$paraText = $titlePara.Range.Text
$hyphenOffset = $paraText.IndexOf("Pseudo-Code")
if ($hyphenOffset -ge 0) {
    $hyphenOffset = $hyphenOffset + 6
    $hyphenPos = $paraStart + $hyphenOffset
    $hyphenRange = $d.Range($hyphenPos, $hyphenPos + 1)
    $hyphenRange.Text = " "
}

# Step 2: Re-split the run so " Pseudo " and "Code" end up as separate runs
# (matching the target layout) without altering any text. Toggling a
# character property on and back off forces the engine to keep the touched
# range as its own run even though the resulting formatting is unchanged.
$paraText = $titlePara.Range.Text
$pseudoOffset = $paraText.IndexOf(" Pseudo Code")
$codeOffset = $paraText.IndexOf("Code")

if ($pseudoOffset -ge 0) {
    $tailRange = $d.Range($paraStart + $pseudoOffset, $paraEnd)
    $tailRange.Bold = 1
    $tailRange.Bold = 0
}

if ($codeOffset -ge 0) {
    $codePos = $paraStart + $codeOffset
    $codeRange = $d.Range($codePos, $codePos + 4)
    $codeRange.Bold = 1
    $codeRange.Bold = 0
}
